$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: rolling quarter-period headers (drop oldest quarter, shift, append newest) ---
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9: rolling publish-date headers (drop oldest date, shift, append newest) ---
$ws.Range("D9").Value = "1400-11-02 (2)"
$ws.Range("E9").Value = "1401-02-25 (12)"
$ws.Range("F9").Value = "1401-05-01 (3)"
$ws.Range("G9").Value = "1401-08-30 (4)"
$ws.Range("H9").Value = "1401-11-19 (3)"
$ws.Range("I9").Value = "1402-02-27 (12)"
$ws.Range("K9").Value = "1401-08-30 (2)"
$ws.Range("L9").Value = "1401-11-19 (2)"
$ws.Range("M9").Value = "1402-02-27 (3)"

# J9 ("1401-05-01") looks like an ISO date, so Excel would auto-convert it to a date
# serial number. Force text format, set the value, then restore the original cell
# formatting (copied from a neighboring cell) so the style index is unaffected.
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "1401-05-01"
$ws.Range("I9").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows: rolling 10-quarter window of financial figures (drop oldest, shift, append newest) ---
# Row 11
$ws.Range("D11").Value = 8060
$ws.Range("E11").Value = 8138
$ws.Range("F11").Value = 11452
$ws.Range("G11").Value = 9867
$ws.Range("H11").Value = 10855
$ws.Range("I11").Value = 10940
$ws.Range("J11").Value = 11624
$ws.Range("K11").Value = 15224
$ws.Range("L11").Value = 18216
$ws.Range("M11").Value = 14720

# Row 12
$ws.Range("D12").Value = -3056
$ws.Range("E12").Value = -3243
$ws.Range("F12").Value = -4595
$ws.Range("G12").Value = -4125
$ws.Range("H12").Value = -4649
$ws.Range("I12").Value = -6272
$ws.Range("J12").Value = -5304
$ws.Range("K12").Value = -5572
$ws.Range("L12").Value = -7040
$ws.Range("M12").Value = -6888

# Row 13
$ws.Range("D13").Value = 5005
$ws.Range("E13").Value = 4895
$ws.Range("F13").Value = 6857
$ws.Range("G13").Value = 5742
$ws.Range("H13").Value = 6206
$ws.Range("I13").Value = 4669
$ws.Range("J13").Value = 6320
$ws.Range("K13").Value = 9652
$ws.Range("L13").Value = 11176
$ws.Range("M13").Value = 7832

# Row 14
$ws.Range("D14").Value = -811
$ws.Range("E14").Value = -1772
$ws.Range("F14").Value = -1244
$ws.Range("G14").Value = -728
$ws.Range("H14").Value = -330
$ws.Range("I14").Value = -576
$ws.Range("J14").Value = -503
$ws.Range("K14").Value = -605
$ws.Range("L14").Value = -2042
$ws.Range("M14").Value = -1676

# Row 16
$ws.Range("D16").Value = 228
$ws.Range("E16").Value = -105
$ws.Range("F16").Value = -8
$ws.Range("G16").Value = -177
$ws.Range("H16").Value = 170
$ws.Range("I16").Value = -102
$ws.Range("J16").Value = 6
$ws.Range("K16").Value = -7
$ws.Range("L16").Value = 244
$ws.Range("M16").Value = 650

# Row 17
$ws.Range("D17").Value = 4422
$ws.Range("E17").Value = 3017
$ws.Range("F17").Value = 5605
$ws.Range("G17").Value = 4837
$ws.Range("H17").Value = 6045
$ws.Range("I17").Value = 3991
$ws.Range("J17").Value = 5823
$ws.Range("K17").Value = 9040
$ws.Range("L17").Value = 9378
$ws.Range("M17").Value = 6806

# Row 18
$ws.Range("D18").Value = -47
$ws.Range("E18").Value = -39
$ws.Range("F18").Value = -58
$ws.Range("G18").Value = -46
$ws.Range("H18").Value = -45
$ws.Range("I18").Value = -37
$ws.Range("J18").Value = -51
$ws.Range("K18").Value = -43
$ws.Range("L18").Value = -41
$ws.Range("M18").Value = -29

# Row 19
$ws.Range("D19").Value = 553
$ws.Range("E19").Value = 84
$ws.Range("F19").Value = 454
$ws.Range("G19").Value = 251
$ws.Range("H19").Value = 620
$ws.Range("I19").Value = 386
$ws.Range("J19").Value = 854
$ws.Range("K19").Value = 794
$ws.Range("L19").Value = 824
$ws.Range("M19").Value = 1027

# Row 20
$ws.Range("D20").Value = 4927
$ws.Range("E20").Value = 3062
$ws.Range("F20").Value = 6002
$ws.Range("G20").Value = 5042
$ws.Range("H20").Value = 6620
$ws.Range("I20").Value = 4340
$ws.Range("J20").Value = 6626
$ws.Range("K20").Value = 9791
$ws.Range("L20").Value = 10161
$ws.Range("M20").Value = 7804

# Row 21
$ws.Range("D21").Value = -713
$ws.Range("E21").Value = 89
$ws.Range("F21").Value = -427
$ws.Range("G21").Value = -1388
$ws.Range("H21").Value = -224
$ws.Range("I21").Value = -596
$ws.Range("J21").Value = -858
$ws.Range("K21").Value = -2160
$ws.Range("L21").Value = -1307
$ws.Range("M21").Value = -858

# Row 22
$ws.Range("D22").Value = 4214
$ws.Range("E22").Value = 3151
$ws.Range("F22").Value = 5575
$ws.Range("G22").Value = 3654
$ws.Range("H22").Value = 6396
$ws.Range("I22").Value = 3744
$ws.Range("J22").Value = 5768
$ws.Range("K22").Value = 7631
$ws.Range("L22").Value = 8853
$ws.Range("M22").Value = 6945

# Row 24
$ws.Range("D24").Value = 4214
$ws.Range("E24").Value = 3151
$ws.Range("F24").Value = 5575
$ws.Range("G24").Value = 3654
$ws.Range("H24").Value = 6396
$ws.Range("I24").Value = 3744
$ws.Range("J24").Value = 5768
$ws.Range("K24").Value = 7631
$ws.Range("L24").Value = 8853
$ws.Range("M24").Value = 6945

# Row 26
$ws.Range("D26").Value = 2392
$ws.Range("E26").Value = 2661
$ws.Range("F26").Value = 2786
$ws.Range("G26").Value = 2490
$ws.Range("H26").Value = 2283
$ws.Range("I26").Value = 2360
$ws.Range("J26").Value = 2212
$ws.Range("K26").Value = 2096
$ws.Range("L26").Value = 1871
$ws.Range("M26").Value = 1430
